$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 and F3 currently store the text "True" in the "Is Active" column.
# Push up a real boolean value instead of the string.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
